# Auto-generated script to apply violent-crime-ytd.xlsx 2024-08-30 data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 5303
$ws.Range('K3').Value = 5474
$ws.Range('F4').Value = 1306
$ws.Range('H4').Value = 1150
$ws.Range('K4').Value = 1134
$ws.Range('K5').Value = 391
$ws.Range('K6').Value = 6077
$ws.Range('F7').Value = 16128
$ws.Range('H7').Value = 16490
$ws.Range('K7').Value = 18379

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('K5').Value = 9
$ws.Range('K6').Value = 18

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K2').Value = 60
$ws.Range('K7').Value = 236

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 370
$ws.Range('K6').Value = 419
$ws.Range('K7').Value = 1234

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 137
$ws.Range('K7').Value = 407

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K6').Value = 230
$ws.Range('K7').Value = 787

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 178
$ws.Range('K3').Value = 205
$ws.Range('K4').Value = 30
$ws.Range('K5').Value = 28
$ws.Range('K6').Value = 181
$ws.Range('K7').Value = 622

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 105
$ws.Range('K7').Value = 417

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 130
$ws.Range('K6').Value = 80
$ws.Range('K7').Value = 312

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K7').Value = 545
$ws.Range('K8').Value = 1234
$ws.Range('K10').Value = 103
$ws.Range('K18').Value = 124
$ws.Range('K19').Value = 541
$ws.Range('K20').Value = 425
$ws.Range('K29').Value = 982
$ws.Range('K31').Value = 201
$ws.Range('K33').Value = 787
$ws.Range('K34').Value = 105
$ws.Range('K36').Value = 243
$ws.Range('K37').Value = 622
$ws.Range('K38').Value = 18
$ws.Range('K41').Value = 128
$ws.Range('K46').Value = 37
$ws.Range('K48').Value = 232
$ws.Range('K49').Value = 102
$ws.Range('K51').Value = 230
$ws.Range('K52').Value = 479
$ws.Range('K53').Value = 236
$ws.Range('K57').Value = 68
$ws.Range('K58').Value = 11
$ws.Range('K60').Value = 115
$ws.Range('F63').Value = 124
$ws.Range('H63').Value = 135
$ws.Range('K63').Value = 50
$ws.Range('K65').Value = 417
$ws.Range('K67').Value = 695
$ws.Range('K71').Value = 58
$ws.Range('K73').Value = 158
$ws.Range('K76').Value = 254
$ws.Range('K77').Value = 129
$ws.Range('K78').Value = 211
$ws.Range('K79').Value = 457
$ws.Range('K83').Value = 407
$ws.Range('K84').Value = 136
$ws.Range('K85').Value = 867
$ws.Range('K89').Value = 268
$ws.Range('K93').Value = 69
$ws.Range('K94').Value = 243
$ws.Range('K96').Value = 199
$ws.Range('K97').Value = 147
$ws.Range('K99').Value = 312
$ws.Range('F101').Value = 16128
$ws.Range('H101').Value = 16490
$ws.Range('K101').Value = 18379

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K3').Value = 49
$ws.Range('K6').Value = 72
$ws.Range('K7').Value = 201

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K3').Value = 248
$ws.Range('K7').Value = 695

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 54
$ws.Range('K7').Value = 136

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K2').Value = 20
$ws.Range('K7').Value = 102

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 281
$ws.Range('K3').Value = 355
$ws.Range('K7').Value = 982

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 32
$ws.Range('K7').Value = 232

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K3').Value = 171
$ws.Range('K6').Value = 172
$ws.Range('K7').Value = 541

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K4').Value = 16
$ws.Range('K7').Value = 254

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K2').Value = 45
$ws.Range('K7').Value = 128

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K2').Value = 29
$ws.Range('K3').Value = 18
$ws.Range('K7').Value = 103

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K6').Value = 76
$ws.Range('K7').Value = 211

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('K3').Value = 10
$ws.Range('K7').Value = 37

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K3').Value = 38
$ws.Range('K7').Value = 199

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 154
$ws.Range('K3').Value = 146
$ws.Range('K4').Value = 30
$ws.Range('K6').Value = 111
$ws.Range('K7').Value = 457

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 140
$ws.Range('K3').Value = 136
$ws.Range('K7').Value = 425

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K3').Value = 39
$ws.Range('K7').Value = 124

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K3').Value = 71
$ws.Range('K7').Value = 243

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K4').Value = 5
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 185
$ws.Range('K3').Value = 175
$ws.Range('K5').Value = 22
$ws.Range('K6').Value = 145
$ws.Range('K7').Value = 545

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K2').Value = 37
$ws.Range('K7').Value = 105

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K2').Value = 69
$ws.Range('K6').Value = 105
$ws.Range('K7').Value = 243

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K4').Value = 11
$ws.Range('K7').Value = 158

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K3').Value = 28
$ws.Range('K7').Value = 147

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 75
$ws.Range('K3').Value = 82
$ws.Range('K6').Value = 80
$ws.Range('K7').Value = 268

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K6').Value = 79
$ws.Range('K7').Value = 230

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K3').Value = 14
$ws.Range('K7').Value = 68

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K3').Value = 34
$ws.Range('K7').Value = 115

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 290
$ws.Range('K3').Value = 293
$ws.Range('K6').Value = 208
$ws.Range('K7').Value = 867

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K5').Value = 1
$ws.Range('K7').Value = 58

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 57
$ws.Range('K7').Value = 129

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 135
$ws.Range('K6').Value = 174
$ws.Range('K7').Value = 479

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range('K6').Value = 9
$ws.Range('K7').Value = 11
